$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.2
$ws.Range("I2").Value = 3.75
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 1.98
$ws.Range("R2").Value = 1.88
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 1.17
$ws.Range("AD2").Value = 9
$ws.Range("AS2").Value = 51
$ws.Range("G3").Value = 2.15
$ws.Range("I3").Value = 4.2
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 5
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1.36
$ws.Range("W3").Value = 6.5
$ws.Range("X3").Value = 1.11
$ws.Range("Y3").Value = 1.67
$ws.Range("Z3").Value = 2.1
$ws.Range("AD3").Value = 8.5
$ws.Range("AS3").Value = 51
$ws.Range("G4").Value = 2.63
$ws.Range("H4").Value = 2.8
$ws.Range("I4").Value = 3.1
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 1.83
$ws.Range("Q4").Value = 2.05
$ws.Range("R4").Value = 1.8
$ws.Range("U4").Value = 4.5
$ws.Range("V4").Value = 1.21
$ws.Range("AA4").Value = 2.2
$ws.Range("AB4").Value = 1.62
$ws.Range("AC4").Value = 6
$ws.Range("AL4").Value = 81
$ws.Range("AN5").Value = 6
$ws.Range("AO5").Value = 13
$ws.Range("G8").Value = 2.4
$ws.Range("I8").Value = 2.9
$ws.Range("J8").Value = 3.1
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("S8").Value = 2.1
$ws.Range("T8").Value = 1.7
$ws.Range("AI8").Value = 9
$ws.Range("AJ8").Value = 6.5
$ws.Range("AO8").Value = 15
$ws.Range("AR8").Value = 26
$ws.Range("G9").Value = 2.25
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 3.4
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 4
$ws.Range("AA9").Value = 1.91
$ws.Range("AB9").Value = 1.8
$ws.Range("AI9").Value = 7.5
$ws.Range("AP9").Value = 13
$ws.Range("H10").Value = 2.92
$ws.Range("I10").Value = 3.4
$ws.Range("K10").Value = 1.98
$ws.Range("L10").Value = 4.1
$ws.Range("N10").Value = 5.7
$ws.Range("O10").Value = 1.47
$ws.Range("P10").Value = 2.52
$ws.Range("S10").Value = 2.35
$ws.Range("T10").Value = 1.53
$ws.Range("W10").Value = 4.1
$ws.Range("X10").Value = 1.19
$ws.Range("Y10").Value = 1.5
$ws.Range("Z10").Value = 2.42
$ws.Range("AD10").Value = 9.5
$ws.Range("AE10").Value = 9
$ws.Range("AF10").Value = 21
$ws.Range("AH10").Value = 35
$ws.Range("AI10").Value = 5.7
$ws.Range("AJ10").Value = 5.8
$ws.Range("AK10").Value = 16.5
$ws.Range("AO10").Value = 17
$ws.Range("AQ10").Value = 50
$ws.Range("AR10").Value = 40
$ws.Range("AS10").Value = 55
$ws.Range("T11").Value = 2
$ws.Range("W11").Value = 3
$ws.Range("X11").Value = 1.36
$ws.Range("J12").Value = 2.87
$ws.Range("K12").Value = 1.87
$ws.Range("M12").Value = 1.07
$ws.Range("O12").Value = 1.47
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 1.88
$ws.Range("X12").Value = 1.11
$ws.Range("AB12").Value = 1.54
$ws.Range("G13").Value = 2.5
$ws.Range("I13").Value = 2.7
$ws.Range("J13").Value = 3.25
$ws.Range("L13").Value = 3.5
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 9
$ws.Range("P13").Value = 3
$ws.Range("S13").Value = 2.1
$ws.Range("T13").Value = 1.7
$ws.Range("X13").Value = 1.22
$ws.Range("AB13").Value = 1.8
$ws.Range("AC13").Value = 8
$ws.Range("AD13").Value = 12
$ws.Range("AF13").Value = 23
$ws.Range("AG13").Value = 21
$ws.Range("AI13").Value = 8.5
$ws.Range("AM13").Value = 301
$ws.Range("AP13").Value = 11
$ws.Range("AQ13").Value = 29
